$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. AP-1 "Privacy" paragraph: replace the WIOA boilerplate with the new
#    client/project authorization sentence. This paragraph has a single run,
#    so a direct Range.Text assignment is safe.
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(7)
$p1.Range.Text = "The Client, as the governing agency of The Project, has authorized the collection of user names and email addresses for the purpose of authenticating to the Project system."

# ---------------------------------------------------------------------------
# 2. AP-2 "Privacy" paragraph: this paragraph currently ends with the run
#    '. Anonymous access is possible, but courses and community participation
#    require an account for which these fields are required:' and is followed
#    by two bulleted list paragraphs and a trailing "Any additional
#    information..." paragraph. We need to:
#      - drop "courses and " from that sentence
#      - fold the two bullets into this same paragraph as plain runs
#        prefixed with a literal "* "
#      - fold the trailing paragraph in too, dropping "as part of coursework
#        or "
#    Directly assigning Range.Text on a sub-range of a multi-run paragraph
#    collapses ALL runs of that paragraph into one - so instead we delete the
#    exact run we want changed (Range.Delete keeps sibling runs intact) and
#    then append the replacement text with InsertAfter, which always creates
#    a fresh run at the insertion point without touching existing runs.
# ---------------------------------------------------------------------------
$bullets = $d.Paragraphs.Item(11)
$oldTailText = ". Anonymous access is possible, but courses and community participation require an account for which these fields are required:"
$newTailText = ". Anonymous access is possible, but community participation require an account for which these fields are required:"

$pStart = $bullets.Range.Start
$paraText = $bullets.Range.Text
$tailIdx = $paraText.IndexOf($oldTailText)
$tailStart = $pStart + $tailIdx
$tailEnd = $tailStart + $oldTailText.Length
$d.Range($tailStart, $tailEnd).Delete()
$bullets.Range.InsertAfter($newTailText)

# ---------------------------------------------------------------------------
# 3. Merge in the former bullet paragraphs (now re-cast as plain runs) and
#    the trailing paragraph. Paragraph marks are removed by deleting the
#    single character at the end of each paragraph's range; the *surviving*
#    paragraph properties are always those of the paragraph that used to
#    follow the deleted mark, so we work from the bottom up - merging
#    paragraph 13 into 14, then 12 into that, then 11 into that - so the
#    final merged paragraph ends up using paragraph 14's clean
#    "FirstParagraph" style (no stray bullet numbering).
# ---------------------------------------------------------------------------

# 3a. Re-word the two bullet paragraphs as "* ..." text (still separate
#     paragraphs at this point, each with a single run, so Range.Text is
#     safe here).
$emailPara = $d.Paragraphs.Item(12)
$emailPara.Range.Text = "* Email address - used for identification."

$namePara = $d.Paragraphs.Item(13)
$namePara.Range.Text = "* First name, last name - used for addressing a logged in user."

# 3b. Fix the trailing paragraph's wording (single run, safe to set Text).
$trailingPara = $d.Paragraphs.Item(14)
$trailingPara.Range.Text = "Any additional information is entered by the user at will to enhance community participation in forums."

# 3c. Merge paragraph 13 forward into paragraph 14, inserting a separating
#     space run first.
$namePara = $d.Paragraphs.Item(13)
$namePara.Range.InsertAfter(" ")
$mark = $d.Range($namePara.Range.End - 1, $namePara.Range.End)
$mark.Delete()

# 3d. Merge paragraph 12 forward into the paragraph produced above.
$emailPara = $d.Paragraphs.Item(12)
$emailPara.Range.InsertAfter(" ")
$mark = $d.Range($emailPara.Range.End - 1, $emailPara.Range.End)
$mark.Delete()

# 3e. Merge paragraph 11 (the bullets intro sentence) forward into the rest.
$bullets = $d.Paragraphs.Item(11)
$bullets.Range.InsertAfter(" ")
$mark = $d.Range($bullets.Range.End - 1, $bullets.Range.End)
$mark.Delete()
